# Result_Analysis.xlsx - "calculated avg transfer time"
#
# Adds, on both the "Saturday Morning" and "Saturday Evening" sheets, a small
# block of computed/typed-in results next to the "Top 5 transfers" table:
#   D24/F24 -> labels "Total Transfers" / "Avg Transfer Wait Time"
#   D25      -> total transfer wait time (minutes), explicit black font
#   F25      -> =D28/D25  (objective / total transfers = avg wait per transfer)
#   D27      -> label "Objective Function"
#   D28      -> objective function value from the AMPL solve
# Also fixes the mislabeled title in A15 of the evening sheet, and leaves the
# selection/active-tab on "Saturday Morning" the way the author ended up
# saving the file.

$wb = $excel.ActiveWorkbook

$wsMorning = $wb.Worksheets.Item("Saturday Morning")
$wsEvening = $wb.Worksheets.Item("Saturday Evening")

# --- Saturday Evening: fix the row-15 title, it used to say "Saturday Morning" ---
$wsEvening.Range("A15").Value = "Saturday Evening"

# --- Objective Function value (entered first, on the evening sheet, then morning) ---
$wsEvening.Range("D27").Value = "Objective Function"
$wsEvening.Range("D28").Value = 117498.25

$wsMorning.Range("D27").Value = "Objective Function"
$wsMorning.Range("D28").Value = 62591.5

# --- Total Transfers / Avg Transfer Wait Time block ---
$wsMorning.Range("D24").Value = "Total Transfers"
$wsMorning.Range("F24").Value = "Avg Transfer Wait Time"
$wsMorning.Range("D25").Value = 47517.25
$wsMorning.Range("D25").Font.Color = 0
$wsMorning.Range("F25").Formula = "=D28/D25"

$wsEvening.Range("D24").Value = "Total Transfers"
$wsEvening.Range("F24").Value = "Avg Transfer Wait Time"
$wsEvening.Range("D25").Value = 43878
$wsEvening.Range("D25").Font.Color = 0
$wsEvening.Range("F25").Formula = "=D28/D25"

# --- Window / selection state: author ended with "Saturday Morning" active,
#     selection at D29, while "Saturday Evening" keeps a D30 selection. ---
[void]$wsEvening.Activate()
[void]$wsEvening.Range("D30").Select()

[void]$wsMorning.Activate()
[void]$wsMorning.Range("D29").Select()
